$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "XQG2DN"
$ws.Range("B9").Value = "Almohadilla + Chip Epson T6716"
$ws.Range("C9").Value = "WF C529 C579 C5210 C5290 C5298 C5299 C5710 C5790 C5799 M5298 M5299 M5799, ET8700, PX S380 S381 S880 S884"
$ws.Range("D9").Value = 30000
$ws.Range("E9").Value = 200000
$ws.Range("F9").Value = 9
$ws.Range("G9").Value = 1
$ws.Range("H9").Formula = "=(E9-D9)*G9"
$ws.Range("I9").Formula = "=D9*F9"
$ws.Range("J9").Value = 270000
